$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Stash the existing bordered-cell style (style index 1) in a scratch cell
# far away from the working area, so we can re-apply it after the old table
# is cleared out.
$ws.Range("A1").Copy()
$ws.Range("Z100").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Remove the old "Search type / Query? / Sources..." lookup table
# (previously A1:C5) entirely.
$ws.Range("A1:C5").Clear()

# --- Re-apply the bordered style onto the new table location, B3:D7.
$ws.Range("Z100").Copy()
$ws.Range("B3:D7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("Z100").Clear()

# --- Fill in the new table contents (shifted down two rows and right one
# column relative to the original table).
$ws.Range("B3").Value = "Search type"
$ws.Range("C3").Value = "Query?"
$ws.Range("D3").Value = "Sources (Netflix, Amazon, Hulu)"

$ws.Range("B4").Value = "Show"
$ws.Range("C4").Value = "No"
$ws.Range("D4").Value = "Yes"

$ws.Range("B5").Value = "Episodes"
$ws.Range("C5").Value = "Id"
$ws.Range("D5").Value = "Yes"

$ws.Range("B6").Value = "Movies"
$ws.Range("C6").Value = "Id"
$ws.Range("D6").Value = "No"

$ws.Range("B7").Value = "Search"
$ws.Range("C7").Value = "Yes"
$ws.Range("D7").Value = "No"

# --- Move the presentation-notes cells from column D to column B.
$ws.Range("D12").Cut($ws.Range("B12"))
$ws.Range("D19").Cut($ws.Range("B19"))
$ws.Range("D20").Cut($ws.Range("B20"))
$ws.Range("D21").Cut($ws.Range("B21"))

# --- Update column widths: column A stays the same; column C's former
# width now belongs to column D, and the new columns B/C get their own
# widths.
$ws.Columns("B").ColumnWidth = 11.6640625
$ws.Columns("C").ColumnWidth = 7.88671875
$ws.Columns("D").ColumnWidth = 28

# --- Update the selected cell to reflect where editing ended up.
$ws.Range("G16").Select()
